$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "Engel Acosta "
$ws.Range("B15").Value = 8298927758
$ws.Range("C15").Value = "16-SIIN-1-095"

$ws.Range("C13").Select()
